$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A77").Value = 19
$ws.Range("B77").Value = "Hibah"
$v = $ws.Range("A77").Value()
Write-Host "Value A77:" $v
